$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = -0.6640704137585214
$ws.Range("J4").Value = 0.457913895639154
$ws.Range("K4").Value = 0.3625523166945925
$ws.Range("L4").Value = 2.7515562354048
